# Apply Flashscore odds updates for 2025-05-22 workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cells {
    param($Sheet, $Row, $Values)
    foreach ($col in $Values.Keys) {
        $Sheet.Range("$col$Row").Value = $Values[$col]
    }
}

# Row 2
Set-Cells $ws 2 @{
    "K"  = 12
    "AF" = 17
}

# Row 10
Set-Cells $ws 10 @{
    "G"  = 1.65
    "H"  = 4.1
    "I"  = 4.5
    "N"  = 1.65
    "O"  = 2.2
    "T"  = 8.5
    "U"  = 9
    "W"  = 13
    "X"  = 12
    "AA" = 8
    "AB" = 15
    "AE" = 15
    "AF" = 26
    "AG" = 15
    "AH" = 51
    "AI" = 34
}

# Row 12
Set-Cells $ws 12 @{
    "N" = 1.62
    "O" = 2.25
}

# Row 13
Set-Cells $ws 13 @{
    "N" = 1.93
    "O" = 1.88
}

# Row 24
Set-Cells $ws 24 @{
    "G"  = 2.7
    "I"  = 2.7
    "J"  = 1.06
    "K"  = 10
    "T"  = 8.5
    "U"  = 13
    "X"  = 23
    "AG" = 10
    "AH" = 26
}

# Row 26
Set-Cells $ws 26 @{
    "G"  = 1.85
    "H"  = 3.75
    "I"  = 3.75
    "W"  = 17
    "AC" = 34
    "AD" = 101
    "AI" = 26
}

# Row 29 - previously empty odds columns (G:AJ), now populated with values
Set-Cells $ws 29 @{
    "G"  = 5.1
    "H"  = 4.15
    "I"  = 1.55
    "J"  = 1.03
    "K"  = 9
    "L"  = 1.19
    "M"  = 4.2
    "N"  = 1.57
    "O"  = 2.25
    "P"  = 1.3
    "Q"  = 3.2
    "R"  = 1.65
    "S"  = 2.1
    "T"  = 17.5
    "U"  = 35
    "V"  = 16
    "W"  = 90
    "X"  = 45
    "Y"  = 40
    "Z"  = 9
    "AA" = 8.25
    "AB" = 14.5
    "AC" = 55
    "AD" = 350
    "AE" = 8.5
    "AF" = 8.25
    "AG" = 8
    "AH" = 11.75
    "AI" = 11.5
    "AJ" = 21
}
